$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BC")

# First batch: new Post-submission related enum/message pairs (rows 7-9)
$batch1 = @(
    @{ Row = 7; Enum = "PcsProject__DuAnDaKetThuc"; Message = "Dự án: {0} đã kết thúc" },
    @{ Row = 8; Enum = "PcsPost__TrangThaiBaiDangKhongHopLe"; Message = "Trạng thái bài đăng không hợp lệ" },
    @{ Row = 9; Enum = "PcsPost__TonTaiBaiDangKhongThuocDuAn"; Message = "Tồn tại bài đăng không thuộc dự án" }
)

foreach ($item in $batch1) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Enum
}
foreach ($item in $batch1) {
    $ws.Cells.Item($item.Row, 2).Value = $item.Message
}

# Second batch: Approve / Check related enum/message pairs (rows 10-11)
$batch2 = @(
    @{ Row = 10; Enum = "PcsPost__DuAnChuaCoBaiDangNaoChuaDuyet"; Message = "Dự án chưa có bài đăng nào chưa duyệt" },
    @{ Row = 11; Enum = "PcsPost__DuAnKhongCoBaiDangNaoHopLe"; Message = "Dự án không có bài đăng nào hợp lệ" }
)

foreach ($item in $batch2) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Enum
}
foreach ($item in $batch2) {
    $ws.Cells.Item($item.Row, 2).Value = $item.Message
}

# Update the selection to match the target state (E11 single cell)
$ws.Range("E11").Select()
